$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

$ws.Range("B2").Value = 0.06743097493015897
$ws.Range("C2").Value = 0.9987681957784403
$ws.Range("D2").Value = 0.206184617309753

$ws.Range("G2").Value = 0.1311458841167526
$ws.Range("H2").Value = 0.991
